$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zuordnungstabelle")

# Fix typos in the "Zuordnungstabelle" sheet
$ws.Range("D6").Value = "Tür blockiert, Auslösung durch Lichtschranken"
$ws.Range("D18").Value = "Lift fährt rauf"
$ws.Range("D19").Value = "Lift fährt runter"

# Page setup (A4 portrait) as left by the author's print-preview/page-setup pass
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection ends up covering the used range with B2 active
$ws.Activate()
[void]$ws.Range("B2:I19").Select()
